$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/3/2024  Through  6/9/2024"

# --- Simple numeric value updates ---
$ws.Range("M14").Value = 50
$ws.Range("M15").Value = 175
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 13
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -5.084745762711
$ws.Range("L16").Value = -24.324324324324
$ws.Range("M16").Value = -52.136752136752
$ws.Range("N16").Value = -85.964912280701
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 69.230769230769
$ws.Range("I17").Value = 127
$ws.Range("J17").Value = 95
$ws.Range("K17").Value = 33.684210526315
$ws.Range("L17").Value = 69.333333333333
$ws.Range("M17").Value = 217.5
$ws.Range("N17").Value = 29.591836734693
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 99
$ws.Range("J18").Value = 122
$ws.Range("K18").Value = -18.852459016393
$ws.Range("L18").Value = 22.222222222222
$ws.Range("M18").Value = -11.607142857142
$ws.Range("N18").Value = -85.462555066079
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 7.142857142857
$ws.Range("I19").Value = 255
$ws.Range("J19").Value = 264
$ws.Range("K19").Value = -3.40909090909
$ws.Range("L19").Value = -8.273381294964
$ws.Range("M19").Value = 18.60465116279
$ws.Range("N19").Value = -0.778210116731
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -69.230769230769
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = -31.428571428571
$ws.Range("I20").Value = 135
$ws.Range("J20").Value = 132
$ws.Range("K20").Value = 2.272727272727
$ws.Range("L20").Value = 23.853211009174
$ws.Range("M20").Value = 46.739130434782
$ws.Range("N20").Value = -93.781667434362
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -13.333333333333
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = -0.813008130081
$ws.Range("I21").Value = 686
$ws.Range("J21").Value = 680
$ws.Range("K21").Value = 0.882352941176
$ws.Range("L21").Value = 9.235668789808
$ws.Range("M21").Value = 17.86941580756
$ws.Range("N21").Value = -81.0706401766
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 5.882352941176
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 21
$ws.Range("K23").Value = -36.363636363636
$ws.Range("L23").Value = -19.230769230769
$ws.Range("M23").Value = 50
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 12.5
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = 5.66037735849
$ws.Range("I24").Value = 533
$ws.Range("J24").Value = 651
$ws.Range("K24").Value = -18.125960061443
$ws.Range("L24").Value = -5.996472663139
$ws.Range("M24").Value = 33.583959899749
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 36.111111111111
$ws.Range("I25").Value = 242
$ws.Range("J25").Value = 236
$ws.Range("K25").Value = 2.542372881355
$ws.Range("L25").Value = 52.201257861635
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 28.571428571428
$ws.Range("I26").Value = 216
$ws.Range("J26").Value = 202
$ws.Range("K26").Value = 6.930693069306
$ws.Range("L26").Value = 10.204081632653
$ws.Range("M26").Value = 21.348314606741
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 0
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 19
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = -9.523809523809
$ws.Range("L28").Value = -36.666666666666
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = -70
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = -70
$ws.Range("G31").Value = 4
$ws.Range("J31").Value = 14
$ws.Range("K31").Value = -92.857142857142

# --- Cells converting from numeric to text (shared-string placeholders) ---
# Pattern: set the literal text value, then copy number-format/style from a
# same-row donor cell that already carries the desired text style (s=14),
# without disturbing the donor cell itself (PasteSpecial copies formats only).
$c = $ws.Range("F14")
$c.Value = "'0"
$ws.Range("E14").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("C15")
$c.Value = "'0"
$ws.Range("D15").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("C22")
$c.Value = "'0"
$ws.Range("N22").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("D23")
$c.Value = "'0"
$ws.Range("N23").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("E23")
$c.Value = "'***.*"
$ws.Range("N23").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("C27")
$c.Value = "'0"
$ws.Range("M27").Copy()
$c.PasteSpecial(-4122)

# --- Cells converting from text back to numeric ---
# Pattern: set the literal numeric value, then copy number-format/style from a
# same-row donor cell that already carries the desired numeric style.
$c = $ws.Range("C23")
$c.Value = 1
$ws.Range("G23").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("D27")
$c.Value = 1
$ws.Range("F27").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("E27")
$c.Value = -100
$ws.Range("H27").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("C28")
$c.Value = 2
$ws.Range("D28").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("D31")
$c.Value = 2
$ws.Range("G31").Copy()
$c.PasteSpecial(-4122)

$c = $ws.Range("E31")
$c.Value = -100
$ws.Range("H31").Copy()
$c.PasteSpecial(-4122)

$excel.CutCopyMode = 0
